$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 613.8077
$ws.Range("I19").Value = 474.1111
$ws.Range("J19").Value = 687.7646999999999
$ws.Range("K19").Value = 474.1111
$ws.Range("L19").Value = 687.7646999999999
$ws.Range("M19").Value = -299.1111
$ws.Range("N19").Value = -1037.7647
$ws.Range("H69").Value = 5735.375
$ws.Range("I69").Value = 6664.3335
$ws.Range("J69").Value = 5178
$ws.Range("K69").Value = 19993.0005
$ws.Range("L69").Value = 15534
$ws.Range("M69").Value = -19119.0005
$ws.Range("N69").Value = -17282
$ws.Range("H72").Value = 5735.375
$ws.Range("I72").Value = 6664.3335
$ws.Range("J72").Value = 5178
$ws.Range("K72").Value = 59979.0015
$ws.Range("L72").Value = 46602
$ws.Range("M72").Value = -55611.0015
$ws.Range("N72").Value = -55338
$ws.Range("H98").Value = 2328.4614
$ws.Range("I98").Value = 1996.6666
$ws.Range("J98").Value = 3075
$ws.Range("K98").Value = 1996.6666
$ws.Range("L98").Value = 3075
$ws.Range("M98").Value = -498.6666
$ws.Range("N98").Value = -6071
$ws.Range("H107").Value = 737
$ws.Range("I107").Value = 577.5
$ws.Range("J107").Value = 1056
$ws.Range("K107").Value = 577.5
$ws.Range("L107").Value = 1056
$ws.Range("M107").Value = 1342.5
$ws.Range("N107").Value = -4896
$ws.Range("H111").Value = 8059
$ws.Range("I111").Value = 4069.7273
$ws.Range("J111").Value = 30000
$ws.Range("K111").Value = 12209.1819
$ws.Range("L111").Value = 90000
$ws.Range("M111").Value = -9142.1819
$ws.Range("N111").Value = -96134
$ws.Range("H112").Value = 2654.5454
$ws.Range("I112").Value = 850
$ws.Range("J112").Value = 3185.2942
$ws.Range("K112").Value = 2550
$ws.Range("L112").Value = 9555.882599999999
$ws.Range("M112").Value = -1442
$ws.Range("N112").Value = -11771.8826
$ws.Range("H113").Value = 4132.609
$ws.Range("I113").Value = 3155
$ws.Range("J113").Value = 4560.3125
$ws.Range("K113").Value = 3155
$ws.Range("L113").Value = 4560.3125
$ws.Range("M113").Value = 99
$ws.Range("N113").Value = -11068.3125
$ws.Range("H115").Value = 1725.4166
$ws.Range("I115").Value = 518.63635
$ws.Range("J115").Value = 15000
$ws.Range("K115").Value = 1555.90905
$ws.Range("L115").Value = 45000
$ws.Range("M115").Value = 11.09095000000002
$ws.Range("N115").Value = -48134
$ws.Range("H116").Value = 4366.375
$ws.Range("I116").Value = 2660.7778
$ws.Range("J116").Value = 5389.7334
$ws.Range("K116").Value = 2660.7778
$ws.Range("L116").Value = 5389.7334
$ws.Range("M116").Value = 781.2222000000002
$ws.Range("N116").Value = -12273.7334
$ws.Range("H122").Value = 2328.4614
$ws.Range("I122").Value = 1996.6666
$ws.Range("J122").Value = 3075
$ws.Range("K122").Value = 5989.9998
$ws.Range("L122").Value = 9225
$ws.Range("M122").Value = -3539.9998
$ws.Range("N122").Value = -14125
$ws.Range("H132").Value = 10395.346
$ws.Range("I132").Value = 8693.277
$ws.Range("K132").Value = 26079.831
$ws.Range("M132").Value = -23549.831
$ws.Range("H135").Value = 1860.3429
$ws.Range("I135").Value = 371.75
$ws.Range("J135").Value = 17738.666
$ws.Range("K135").Value = 3345.75
$ws.Range("L135").Value = 159647.994
$ws.Range("M135").Value = -810.75
$ws.Range("N135").Value = -164717.994

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2457
$ws.Range("I2").Value = 1512.7333
$ws.Range("K2").Value = 1512.7333
$ws.Range("M2").Value = -1399.7333
$ws.Range("H32").Value = 7759.6113
$ws.Range("I32").Value = 8046.4243
$ws.Range("J32").Value = 4604.6665
$ws.Range("K32").Value = 8046.4243
$ws.Range("L32").Value = 4604.6665
$ws.Range("M32").Value = -7759.4243
$ws.Range("N32").Value = -5178.6665
$ws.Range("H45").Value = 3455.087
$ws.Range("I45").Value = 2374.2856
$ws.Range("J45").Value = 5136.3335
$ws.Range("K45").Value = 2374.2856
$ws.Range("L45").Value = 5136.3335
$ws.Range("M45").Value = -1997.2856
$ws.Range("N45").Value = -5890.3335
$ws.Range("H88").Value = 2427.5715
$ws.Range("I88").Value = 2197.2
$ws.Range("J88").Value = 3003.5
$ws.Range("K88").Value = 2197.2
$ws.Range("L88").Value = 3003.5
$ws.Range("M88").Value = -1791.2
$ws.Range("N88").Value = -3815.5
$ws.Range("H91").Value = 2427.5715
$ws.Range("I91").Value = 2197.2
$ws.Range("J91").Value = 3003.5
$ws.Range("K91").Value = 2197.2
$ws.Range("L91").Value = 3003.5
$ws.Range("M91").Value = -793.1999999999998
$ws.Range("N91").Value = -5811.5
$ws.Range("H96").Value = 24780
$ws.Range("J96").Value = 24780
$ws.Range("L96").Value = 24780
$ws.Range("N96").Value = -30272
$ws.Range("H110").Value = 5687.2856
$ws.Range("I110").Value = 3202.75
$ws.Range("K110").Value = 3202.75
$ws.Range("M110").Value = -1157.75
$ws.Range("H116").Value = 2457
$ws.Range("I116").Value = 1512.7333
$ws.Range("K116").Value = 1512.7333
$ws.Range("M116").Value = 781.2666999999999
$ws.Range("H132").Value = 4780.1
$ws.Range("I132").Value = 1959.7333
$ws.Range("J132").Value = 13241.2
$ws.Range("K132").Value = 5879.199900000001
$ws.Range("L132").Value = 39723.60000000001
$ws.Range("M132").Value = -3349.199900000001
$ws.Range("N132").Value = -44783.60000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2457
$ws.Range("I3").Value = 1512.7333
$ws.Range("K3").Value = 1512.7333
$ws.Range("M3").Value = -1398.7333
$ws.Range("H86").Value = 1464.5
$ws.Range("I86").Value = 1490.8334
$ws.Range("J86").Value = 1425
$ws.Range("K86").Value = 1490.8334
$ws.Range("L86").Value = 1425
$ws.Range("M86").Value = -367.8334
$ws.Range("N86").Value = -3671
$ws.Range("H89").Value = 1464.5
$ws.Range("I89").Value = 1490.8334
$ws.Range("J89").Value = 1425
$ws.Range("K89").Value = 7454.166999999999
$ws.Range("L89").Value = 7125
$ws.Range("M89").Value = -1838.166999999999
$ws.Range("N89").Value = -18357
$ws.Range("H107").Value = 1951.9166
$ws.Range("I107").Value = 1856.6364
$ws.Range("K107").Value = 1856.6364
$ws.Range("M107").Value = 63.36359999999991
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4833.5
$ws.Range("I16").Value = 6155.5
$ws.Range("J16").Value = 3952.1667
$ws.Range("K16").Value = 6155.5
$ws.Range("L16").Value = 3952.1667
$ws.Range("M16").Value = -5868.5
$ws.Range("N16").Value = -4526.1667
$ws.Range("H31").Value = 1840.75
$ws.Range("I31").Value = 1400.3959
$ws.Range("K31").Value = 1400.3959
$ws.Range("M31").Value = -1105.3959
$ws.Range("H34").Value = 1840.75
$ws.Range("I34").Value = 1400.3959
$ws.Range("K34").Value = 1400.3959
$ws.Range("M34").Value = -1198.3959
$ws.Range("H62").Value = 9583.941000000001
$ws.Range("I62").Value = 2559.6667
$ws.Range("J62").Value = 17486.25
$ws.Range("K62").Value = 2559.6667
$ws.Range("L62").Value = 17486.25
$ws.Range("M62").Value = -1935.6667
$ws.Range("N62").Value = -18734.25
$ws.Range("H65").Value = 9583.941000000001
$ws.Range("I65").Value = 2559.6667
$ws.Range("J65").Value = 17486.25
$ws.Range("K65").Value = 12798.3335
$ws.Range("L65").Value = 87431.25
$ws.Range("M65").Value = -9678.333500000001
$ws.Range("N65").Value = -93671.25
$ws.Range("H113").Value = 4833.5
$ws.Range("I113").Value = 6155.5
$ws.Range("J113").Value = 3952.1667
$ws.Range("K113").Value = 6155.5
$ws.Range("L113").Value = 3952.1667
$ws.Range("M113").Value = -3985.5
$ws.Range("N113").Value = -8292.1667

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2298.4546
$ws.Range("I5").Value = 552
$ws.Range("K5").Value = 1656
$ws.Range("M5").Value = -1544
$ws.Range("H113").Value = 523.3214
$ws.Range("J113").Value = 561.8148
$ws.Range("L113").Value = 1685.4444
$ws.Range("N113").Value = -6025.4444
$ws.Range("H122").Value = 3458.552
$ws.Range("J122").Value = 3776.5518
$ws.Range("L122").Value = 33988.9662
$ws.Range("N122").Value = -38888.9662
$ws.Range("H131").Value = 946.90625
$ws.Range("J131").Value = 1054.2084
$ws.Range("L131").Value = 3162.6252
$ws.Range("N131").Value = -13242.6252
$ws.Range("H135").Value = 2298.4546
$ws.Range("I135").Value = 552
$ws.Range("K135").Value = 4968
$ws.Range("M135").Value = -2433

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6882.8667
$ws.Range("I70").Value = 5262
$ws.Range("J70").Value = 8735.286
$ws.Range("K70").Value = 5262
$ws.Range("L70").Value = 8735.286
$ws.Range("M70").Value = -4992
$ws.Range("N70").Value = -9275.286
$ws.Range("H73").Value = 6882.8667
$ws.Range("I73").Value = 5262
$ws.Range("J73").Value = 8735.286
$ws.Range("K73").Value = 5262
$ws.Range("L73").Value = 8735.286
$ws.Range("M73").Value = -4326
$ws.Range("N73").Value = -10607.286
$ws.Range("H113").Value = 2130.1428
$ws.Range("I113").Value = 1818.5
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 1818.5
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 351.5
$ws.Range("N113").Value = -8340

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 47621540
$ws.Range("I40").Value = 58825240
$ws.Range("J40").Value = 5800
$ws.Range("K40").Value = 58825240
$ws.Range("L40").Value = 5800
$ws.Range("M40").Value = -58825104
$ws.Range("N40").Value = -6072
$ws.Range("H61").Value = 4272.684
$ws.Range("I61").Value = 3443.7273
$ws.Range("K61").Value = 3443.7273
$ws.Range("M61").Value = -3241.7273
$ws.Range("H113").Value = 4272.684
$ws.Range("I113").Value = 3443.7273
$ws.Range("K113").Value = 3443.7273
$ws.Range("M113").Value = -1273.7273

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3575
$ws.Range("I81").Value = 10001
$ws.Range("J81").Value = 1433
$ws.Range("K81").Value = 20002
$ws.Range("L81").Value = 2866
$ws.Range("M81").Value = -18941
$ws.Range("N81").Value = -4988
$ws.Range("H84").Value = 3575
$ws.Range("I84").Value = 10001
$ws.Range("J84").Value = 1433
$ws.Range("K84").Value = 100010
$ws.Range("L84").Value = 14330
$ws.Range("M84").Value = -94706
$ws.Range("N84").Value = -24938
$ws.Range("H126").Value = 4918.8237
$ws.Range("I126").Value = 5834.037
$ws.Range("J126").Value = 1388.7142
$ws.Range("K126").Value = 17502.111
$ws.Range("L126").Value = 4166.142599999999
$ws.Range("M126").Value = -15032.111
$ws.Range("N126").Value = -9106.142599999999
